$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("H6").Value = 76
$ws.Range("L6").Value = "PRUEBAS"

# Row 7
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 34
$ws.Range("I7").Value = 10

# Row 8
$ws.Range("H8").Value = 434

# Row 9
$ws.Range("E9").Value = 43434
$ws.Range("G9").Value = 434

# Row 10
$ws.Range("C10").Value = 45
$ws.Range("E10").Value = 565
$ws.Range("H10").Value = 434

# Row 12
$ws.Range("I12").Value = 20

# Row 13
$ws.Range("C13").Value = 4343
$ws.Range("G13").Value = 434

# Row 14 (added before row 11's GAM so the shared-string table order matches)
$ws.Range("L14").Value = "SADÑLNSA"

# Row 11
$ws.Range("L11").Value = "GAM"

# Row 15
$ws.Range("E15").Value = 777
$ws.Range("H15").Value = 434

# Row 17
$ws.Range("D17").Value = 43
$ws.Range("I17").Value = 3

# Row 18
$ws.Range("D18").Value = 24

# Move the active selection to C19, as reflected in the saved view state
$ws.Range("C19").Select()
